$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 369.125
$ws.Range("J33").Value = 376
$ws.Range("L33").Value = 376
$ws.Range("N33").Value = -834
$ws.Range("H57").Value = 31854
$ws.Range("J57").Value = 31854
$ws.Range("L57").Value = 95562
$ws.Range("N57").Value = -96560
$ws.Range("H98").Value = 4822.6665
$ws.Range("I98").Value = 5262.5
$ws.Range("J98").Value = 1304
$ws.Range("K98").Value = 5262.5
$ws.Range("L98").Value = 1304
$ws.Range("M98").Value = -3764.5
$ws.Range("N98").Value = -4300
$ws.Range("H122").Value = 4822.6665
$ws.Range("I122").Value = 5262.5
$ws.Range("J122").Value = 1304
$ws.Range("K122").Value = 15787.5
$ws.Range("L122").Value = 3912
$ws.Range("M122").Value = -13337.5
$ws.Range("N122").Value = -8812
$ws.Range("H129").Value = 19608646
$ws.Range("I129").Value = 616.7143
$ws.Range("J129").Value = 111112780
$ws.Range("K129").Value = 1850.1429
$ws.Range("L129").Value = 333338340
$ws.Range("M129").Value = 3149.8571
$ws.Range("N129").Value = -333348340
$ws.Range("H132").Value = 444531.6
$ws.Range("J132").Value = 14311.75
$ws.Range("L132").Value = 42935.25
$ws.Range("N132").Value = -47995.25
$ws.Range("H135").Value = 3302.459
$ws.Range("I135").Value = 1444.3773
$ws.Range("K135").Value = 12999.3957
$ws.Range("M135").Value = -10464.3957
$ws.Range("H137").Value = 9739.444
$ws.Range("I137").Value = 5836.647
$ws.Range("J137").Value = 16374.2
$ws.Range("K137").Value = 17509.941
$ws.Range("L137").Value = 49122.60000000001
$ws.Range("M137").Value = -14959.941
$ws.Range("N137").Value = -54222.60000000001
$ws.Range("H141").Value = 886.46155
$ws.Range("I141").Value = 793.65216
$ws.Range("J141").Value = 1598
$ws.Range("K141").Value = 2380.95648
$ws.Range("L141").Value = 4794
$ws.Range("M141").Value = 2799.04352
$ws.Range("N141").Value = -15154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2651248.8
$ws.Range("I32").Value = 5205.5356
$ws.Range("K32").Value = 5205.5356
$ws.Range("M32").Value = -4918.5356
$ws.Range("H45").Value = 4100
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2623
$ws.Range("H61").Value = 5460.6
$ws.Range("I61").Value = 6071.533
$ws.Range("K61").Value = 6071.533
$ws.Range("M61").Value = -5859.533
$ws.Range("H74").Value = 4378.5884
$ws.Range("I74").Value = 5274.0835
$ws.Range("K74").Value = 5274.0835
$ws.Range("M74").Value = -4400.0835
$ws.Range("H77").Value = 4378.5884
$ws.Range("I77").Value = 5274.0835
$ws.Range("K77").Value = 26370.4175
$ws.Range("M77").Value = -22002.4175
$ws.Range("H122").Value = 1926.5869
$ws.Range("I122").Value = 1432.1765
$ws.Range("K122").Value = 4296.529500000001
$ws.Range("M122").Value = -1846.529500000001
$ws.Range("H136").Value = 5460.6
$ws.Range("I136").Value = 6071.533
$ws.Range("K136").Value = 18214.599
$ws.Range("M136").Value = -15664.599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 41689.8
$ws.Range("J74").Value = 41689.8
$ws.Range("L74").Value = 41689.8
$ws.Range("N74").Value = -43561.8
$ws.Range("H77").Value = 41689.8
$ws.Range("J77").Value = 41689.8
$ws.Range("L77").Value = 125069.4
$ws.Range("N77").Value = -134429.4
$ws.Range("H134").Value = 2182722.8
$ws.Range("I134").Value = 3133320.5
$ws.Range("J134").Value = 9928.143
$ws.Range("K134").Value = 9399961.5
$ws.Range("L134").Value = 29784.429
$ws.Range("M134").Value = -9397426.5
$ws.Range("N134").Value = -34854.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4787.364
$ws.Range("I31").Value = 2313.4167
$ws.Range("J31").Value = 6201.048
$ws.Range("K31").Value = 2313.4167
$ws.Range("L31").Value = 6201.048
$ws.Range("M31").Value = -2018.4167
$ws.Range("N31").Value = -6791.048
$ws.Range("H34").Value = 4787.364
$ws.Range("I34").Value = 2313.4167
$ws.Range("J34").Value = 6201.048
$ws.Range("K34").Value = 2313.4167
$ws.Range("L34").Value = 6201.048
$ws.Range("M34").Value = -2111.4167
$ws.Range("N34").Value = -6605.048
$ws.Range("H58").Value = 22733134
$ws.Range("I58").Value = 30307368
$ws.Range("K58").Value = 30307368
$ws.Range("M58").Value = -30307165
$ws.Range("H132").Value = 4228.5
$ws.Range("I132").Value = 3936.5312
$ws.Range("J132").Value = 8900
$ws.Range("K132").Value = 11809.5936
$ws.Range("L132").Value = 26700
$ws.Range("M132").Value = -9279.5936
$ws.Range("N132").Value = -31760
$ws.Range("H134").Value = 33339096
$ws.Range("I134").Value = 50005092
$ws.Range("J134").Value = 7104.1
$ws.Range("K134").Value = 150015276
$ws.Range("L134").Value = 21312.3
$ws.Range("M134").Value = -150012741
$ws.Range("N134").Value = -26382.3
$ws.Range("H136").Value = 22733134
$ws.Range("I136").Value = 30307368
$ws.Range("K136").Value = 90922104
$ws.Range("M136").Value = -90919554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 750
$ws.Range("I51").Value = 750
$ws.Range("K51").Value = 2250
$ws.Range("M51").Value = -1790
$ws.Range("H93").Value = 11429.2
$ws.Range("I93").Value = 8036.5
$ws.Range("K93").Value = 24109.5
$ws.Range("M93").Value = -22237.5
$ws.Range("H112").Value = 10110.571
$ws.Range("I112").Value = 7274.8
$ws.Range("K112").Value = 21824.4
$ws.Range("M112").Value = -20716.4
$ws.Range("H131").Value = 37685360
$ws.Range("I131").Value = 53338916
$ws.Range("K131").Value = 160016748
$ws.Range("M131").Value = -160011708
$ws.Range("H137").Value = 2264.889
$ws.Range("I137").Value = 2134.6667
$ws.Range("J137").Value = 2290.9333
$ws.Range("K137").Value = 6404.000100000001
$ws.Range("L137").Value = 6872.7999
$ws.Range("M137").Value = -1304.000100000001
$ws.Range("N137").Value = -17072.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17547174
$ws.Range("I132").Value = 24393346
$ws.Range("J132").Value = 3855.375
$ws.Range("K132").Value = 73180038
$ws.Range("L132").Value = 11566.125
$ws.Range("M132").Value = -73177508
$ws.Range("N132").Value = -16626.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 48391976
$ws.Range("I136").Value = 19235816
$ws.Range("J136").Value = 200004020
$ws.Range("K136").Value = 57707448
$ws.Range("L136").Value = 600012060
$ws.Range("M136").Value = -57704898
$ws.Range("N136").Value = -600017160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 21666.5
$ws.Range("H82").Value = 466665.66
$ws.Range("J82").Value = 466665.66
$ws.Range("L82").Value = 466665.66
$ws.Range("N82").Value = -467431.66
$ws.Range("H85").Value = 466665.66
$ws.Range("J85").Value = 466665.66
$ws.Range("L85").Value = 466665.66
$ws.Range("N85").Value = -469317.66
$ws.Range("H113").Value = 9805609
$ws.Range("I113").Value = 13890492
$ws.Range("K113").Value = 41671476
$ws.Range("M113").Value = -41669306
$ws.Range("H132").Value = 3500.0986
$ws.Range("I132").Value = 2919.2593
$ws.Range("J132").Value = 5345.1177
$ws.Range("K132").Value = 8757.777900000001
$ws.Range("L132").Value = 16035.3531
$ws.Range("M132").Value = -6227.777900000001
$ws.Range("N132").Value = -21095.3531
$ws.Range("H136").Value = 12832053
$ws.Range("I136").Value = 15632659
$ws.Range("K136").Value = 46897977
$ws.Range("M136").Value = -46895427
